$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Prepare rows 144-146 by copying formatting (and placeholder values) from row 143 ---
$ws.Range("A143:V143").Copy($ws.Range("A144:V144"))
$ws.Range("A143:V143").Copy($ws.Range("A145:V145"))
$ws.Range("A143:V143").Copy($ws.Range("A146:V146"))

# --- Swap mis-ordered match rows ---

# Row 5
$ws.Range("F5").Value = "FC Porto B"
$ws.Range("H5").Value = "Tondela"
$ws.Range("J5").Value = 1.76
$ws.Range("L5").Value = 2.4
$ws.Range("M5").Value = "13/08/2023 11:59"
$ws.Range("N5").Value = 3.72
$ws.Range("P5").Value = 3.33
$ws.Range("Q5").Value = "13/08/2023 11:55"
$ws.Range("R5").Value = 4.79
$ws.Range("T5").Value = 3.16
$ws.Range("U5").Value = "13/08/2023 11:59"
$ws.Range("V5").Value = "https://www.betexplorer.com/football/portugal/liga-portugal-2/fc-porto-tondela/MejJgCPJ/"

# Row 6
$ws.Range("F6").Value = "Academico Viseu"
$ws.Range("H6").Value = "Vilaverdense"
$ws.Range("J6").Value = 1.81
$ws.Range("L6").Value = 1.75
$ws.Range("M6").Value = "13/08/2023 11:52"
$ws.Range("N6").Value = 3.65
$ws.Range("P6").Value = 3.88
$ws.Range("Q6").Value = "13/08/2023 11:52"
$ws.Range("R6").Value = 4.53
$ws.Range("T6").Value = 4.8
$ws.Range("U6").Value = "13/08/2023 11:52"
$ws.Range("V6").Value = "https://www.betexplorer.com/football/portugal/liga-portugal-2/academico-viseu-vilaverdense-fc/Yeoap8n6/"

# Row 39
$ws.Range("F39").Value = "FC Porto B"
$ws.Range("G39").Value = 0
$ws.Range("H39").Value = "Maritimo"
$ws.Range("I39").Value = 2
$ws.Range("J39").Value = 2.76
$ws.Range("K39").Value = "13/09/2023 21:12"
$ws.Range("L39").Value = 2.49
$ws.Range("M39").Value = "16/09/2023 11:00"
$ws.Range("N39").Value = 3.22
$ws.Range("O39").Value = "13/09/2023 21:12"
$ws.Range("P39").Value = 3.34
$ws.Range("Q39").Value = "16/09/2023 10:30"
$ws.Range("R39").Value = 2.58
$ws.Range("S39").Value = "13/09/2023 21:12"
$ws.Range("T39").Value = 3.01
$ws.Range("U39").Value = "16/09/2023 11:00"
$ws.Range("V39").Value = "https://www.betexplorer.com/football/portugal/liga-portugal-2/fc-porto-maritimo/tORXnMP1/"

# Row 40
$ws.Range("F40").Value = "AVS"
$ws.Range("G40").Value = 2
$ws.Range("H40").Value = "Vilaverdense"
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 1.75
$ws.Range("K40").Value = "12/09/2023 12:12"
$ws.Range("L40").Value = 1.82
$ws.Range("M40").Value = "16/09/2023 11:53"
$ws.Range("N40").Value = 3.71
$ws.Range("O40").Value = "12/09/2023 12:12"
$ws.Range("P40").Value = 3.69
$ws.Range("Q40").Value = "16/09/2023 11:53"
$ws.Range("R40").Value = 4.89
$ws.Range("S40").Value = "12/09/2023 12:12"
$ws.Range("T40").Value = 4.66
$ws.Range("U40").Value = "16/09/2023 11:53"
$ws.Range("V40").Value = "https://www.betexplorer.com/football/portugal/liga-portugal-2/avs-vilaverdense-fc/rVCSkOuq/"

# Row 74
$ws.Range("F74").Value = "FC Porto B"
$ws.Range("G74").Value = 2
$ws.Range("H74").Value = "Feirense"
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 1.98
$ws.Range("L74").Value = 1.84
$ws.Range("M74").Value = "04/11/2023 11:59"
$ws.Range("N74").Value = 3.57
$ws.Range("P74").Value = 3.78
$ws.Range("Q74").Value = "04/11/2023 11:59"
$ws.Range("R74").Value = 3.87
$ws.Range("T74").Value = 4.41
$ws.Range("U74").Value = "04/11/2023 11:58"
$ws.Range("V74").Value = "https://www.betexplorer.com/football/portugal/liga-portugal-2/fc-porto-feirense/jTL6QSDN/"

# Row 75
$ws.Range("F75").Value = "Mafra"
$ws.Range("G75").Value = 0
$ws.Range("H75").Value = "Leixoes"
$ws.Range("I75").Value = 1
$ws.Range("J75").Value = 1.88
$ws.Range("L75").Value = 1.93
$ws.Range("M75").Value = "04/11/2023 11:48"
$ws.Range("N75").Value = 3.73
$ws.Range("P75").Value = 3.55
$ws.Range("Q75").Value = "04/11/2023 11:51"
$ws.Range("R75").Value = 3.8
$ws.Range("T75").Value = 4.21
$ws.Range("U75").Value = "04/11/2023 11:51"
$ws.Range("V75").Value = "https://www.betexplorer.com/football/portugal/liga-portugal-2/mafra-leixoes/YiBBPnTT/"

# Row 77
$ws.Range("F77").Value = "Nacional"
$ws.Range("H77").Value = "Santa Clara"
$ws.Range("I77").Value = 1
$ws.Range("J77").Value = 2.98
$ws.Range("K77").Value = "01/11/2023 16:12"
$ws.Range("L77").Value = 2.81
$ws.Range("M77").Value = "04/11/2023 18:58"
$ws.Range("N77").Value = 3.27
$ws.Range("O77").Value = "01/11/2023 16:12"
$ws.Range("P77").Value = 3.23
$ws.Range("Q77").Value = "04/11/2023 18:52"
$ws.Range("R77").Value = 2.39
$ws.Range("S77").Value = "01/11/2023 16:12"
$ws.Range("T77").Value = 2.72
$ws.Range("U77").Value = "04/11/2023 18:52"
$ws.Range("V77").Value = "https://www.betexplorer.com/football/portugal/liga-portugal-2/nacional-santa-clara/xQH2R8bH/"

# Row 78
$ws.Range("F78").Value = "Benfica B"
$ws.Range("H78").Value = "Penafiel"
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 1.91
$ws.Range("K78").Value = "29/10/2023 16:42"
$ws.Range("L78").Value = 2.26
$ws.Range("M78").Value = "04/11/2023 18:53"
$ws.Range("N78").Value = 3.67
$ws.Range("O78").Value = "29/10/2023 16:42"
$ws.Range("P78").Value = 3.59
$ws.Range("Q78").Value = "04/11/2023 18:53"
$ws.Range("R78").Value = 4.01
$ws.Range("S78").Value = "29/10/2023 16:42"
$ws.Range("T78").Value = 3.2
$ws.Range("U78").Value = "04/11/2023 18:53"
$ws.Range("V78").Value = "https://www.betexplorer.com/football/portugal/liga-portugal-2/benfica-penafiel/xjmbUAEb/"

# Row 111
$ws.Range("F111").Value = "Pacos Ferreira"
$ws.Range("H111").Value = "Penafiel"
$ws.Range("I111").Value = 1
$ws.Range("J111").Value = 1.81
$ws.Range("K111").Value = "03/12/2023 15:12"
$ws.Range("L111").Value = 1.89
$ws.Range("M111").Value = "09/12/2023 11:58"
$ws.Range("N111").Value = 3.53
$ws.Range("O111").Value = "03/12/2023 15:12"
$ws.Range("P111").Value = 3.42
$ws.Range("Q111").Value = "09/12/2023 11:58"
$ws.Range("R111").Value = 4.33
$ws.Range("S111").Value = "03/12/2023 15:12"
$ws.Range("T111").Value = 4.65
$ws.Range("U111").Value = "09/12/2023 11:58"
$ws.Range("V111").Value = "https://www.betexplorer.com/football/portugal/liga-portugal-2/pacos-ferreira-penafiel/nyiwnzt1/"

# Row 112
$ws.Range("F112").Value = "FC Porto B"
$ws.Range("H112").Value = "AVS"
$ws.Range("I112").Value = 3
$ws.Range("J112").Value = 2.36
$ws.Range("K112").Value = "02/12/2023 21:42"
$ws.Range("L112").Value = 2.13
$ws.Range("M112").Value = "09/12/2023 11:57"
$ws.Range("N112").Value = 3.26
$ws.Range("O112").Value = "02/12/2023 21:42"
$ws.Range("P112").Value = 3.4
$ws.Range("Q112").Value = "09/12/2023 11:52"
$ws.Range("R112").Value = 3
$ws.Range("S112").Value = "02/12/2023 21:42"
$ws.Range("T112").Value = 3.7
$ws.Range("U112").Value = "09/12/2023 11:57"
$ws.Range("V112").Value = "https://www.betexplorer.com/football/portugal/liga-portugal-2/fc-porto-avs/hfgU6yIQ/"

# Row 134
$ws.Range("F134").Value = "Pacos Ferreira"
$ws.Range("G134").Value = 2
$ws.Range("H134").Value = "Benfica B"
$ws.Range("I134").Value = 2
$ws.Range("J134").Value = 2.02
$ws.Range("K134").Value = "23/12/2023 15:12"
$ws.Range("L134").Value = 2.19
$ws.Range("M134").Value = "30/12/2023 15:00"
$ws.Range("N134").Value = 3.41
$ws.Range("O134").Value = "23/12/2023 15:12"
$ws.Range("P134").Value = 3.21
$ws.Range("Q134").Value = "30/12/2023 15:00"
$ws.Range("R134").Value = 3.59
$ws.Range("S134").Value = "23/12/2023 15:12"
$ws.Range("T134").Value = 3.76
$ws.Range("U134").Value = "30/12/2023 15:00"
$ws.Range("V134").Value = "https://www.betexplorer.com/football/portugal/liga-portugal-2/pacos-ferreira-benfica/tdLEC05n/"

# Row 135
$ws.Range("F135").Value = "Tondela"
$ws.Range("G135").Value = 1
$ws.Range("H135").Value = "Oliveirense"
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 1.5
$ws.Range("K135").Value = "23/12/2023 12:12"
$ws.Range("L135").Value = 1.68
$ws.Range("M135").Value = "30/12/2023 14:51"
$ws.Range("N135").Value = 4.22
$ws.Range("O135").Value = "23/12/2023 12:12"
$ws.Range("P135").Value = 3.94
$ws.Range("Q135").Value = "30/12/2023 14:59"
$ws.Range("R135").Value = 6.04
$ws.Range("S135").Value = "23/12/2023 12:12"
$ws.Range("T135").Value = 5.31
$ws.Range("U135").Value = "30/12/2023 14:59"
$ws.Range("V135").Value = "https://www.betexplorer.com/football/portugal/liga-portugal-2/tondela-oliveirense/ADUrIMsP/"

# --- Set values for the 3 newly appended rows ---

# Row 144
$ws.Range("A144").Value = 143
$ws.Range("B144").Value = "portugal"
$ws.Range("C144").Value = "liga-portugal-2"
$ws.Range("D144").Value = "2023-2024"
$ws.Range("E144").Value = 45298.5
$ws.Range("F144").Value = "AVS"
$ws.Range("G144").Value = 1
$ws.Range("H144").Value = "Leixoes"
$ws.Range("I144").Value = 3
$ws.Range("J144").Value = 1.76
$ws.Range("K144").Value = "31/12/2024 12:12"
$ws.Range("L144").Value = 1.71
$ws.Range("M144").Value = "07/01/2024 11:42"
$ws.Range("N144").Value = 3.57
$ws.Range("O144").Value = "31/12/2024 12:12"
$ws.Range("P144").Value = 3.61
$ws.Range("Q144").Value = "07/01/2024 11:42"
$ws.Range("R144").Value = 4.54
$ws.Range("S144").Value = "31/12/2024 12:12"
$ws.Range("T144").Value = 5.69
$ws.Range("U144").Value = "07/01/2024 11:42"
$ws.Range("V144").Value = "https://www.betexplorer.com/football/portugal/liga-portugal-2/avs-leixoes/GEDR9bk5/"

# Row 145
$ws.Range("A145").Value = 144
$ws.Range("B145").Value = "portugal"
$ws.Range("C145").Value = "liga-portugal-2"
$ws.Range("D145").Value = "2023-2024"
$ws.Range("E145").Value = 45298.625
$ws.Range("F145").Value = "Pacos Ferreira"
$ws.Range("G145").Value = 3
$ws.Range("H145").Value = "FC Porto B"
$ws.Range("I145").Value = 0
$ws.Range("J145").Value = 2.78
$ws.Range("K145").Value = "31/12/2024 15:12"
$ws.Range("L145").Value = 2.72
$ws.Range("M145").Value = "07/01/2024 14:52"
$ws.Range("N145").Value = 3.46
$ws.Range("O145").Value = "31/12/2024 15:12"
$ws.Range("P145").Value = 3.38
$ws.Range("Q145").Value = "07/01/2024 14:51"
$ws.Range("R145").Value = 2.43
$ws.Range("S145").Value = "31/12/2024 15:12"
$ws.Range("T145").Value = 2.71
$ws.Range("U145").Value = "07/01/2024 14:52"
$ws.Range("V145").Value = "https://www.betexplorer.com/football/portugal/liga-portugal-2/pacos-ferreira-fc-porto/xOENAvza/"

# Row 146
$ws.Range("A146").Value = 145
$ws.Range("B146").Value = "portugal"
$ws.Range("C146").Value = "liga-portugal-2"
$ws.Range("D146").Value = "2023-2024"
$ws.Range("E146").Value = 45298.6875
$ws.Range("F146").Value = "Feirense"
$ws.Range("G146").Value = 3
$ws.Range("H146").Value = "Torreense"
$ws.Range("I146").Value = 1
$ws.Range("J146").Value = 2.73
$ws.Range("K146").Value = "31/12/2024 15:12"
$ws.Range("L146").Value = 3.03
$ws.Range("M146").Value = "07/01/2024 16:27"
$ws.Range("N146").Value = 3.05
$ws.Range("O146").Value = "31/12/2024 15:12"
$ws.Range("P146").Value = 3.12
$ws.Range("Q146").Value = "07/01/2024 16:27"
$ws.Range("R146").Value = 2.71
$ws.Range("S146").Value = "31/12/2024 15:12"
$ws.Range("T146").Value = 2.61
$ws.Range("U146").Value = "07/01/2024 16:27"
$ws.Range("V146").Value = "https://www.betexplorer.com/football/portugal/liga-portugal-2/feirense-torreense/AyDV8I4B/"
